# Apply cryptos.xlsx price/volume refresh described by the commit diff.
# Each data row (2-51) on the active sheet holds: B=Coin, C=Link, D=Price, E=Volume(1h).
# D/E are plain text cells (not numbers) in the source file, so numeric-looking
# Price strings are written with a leading apostrophe to force text, then
# ClearFormats() strips the transient "quote prefix" cell format Excel applies,
# leaving the cell with its original (default) style and a plain text value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '43.422.29'
$ws.Range('E2').Value = '  +2.86%  '

# Row 3
$ws.Range('D3').Value = '2.309.40'
$ws.Range('E3').Value = '  +1.84%  '

# Row 4
$ws.Range('E4').Value = '  -0.02%  '

# Row 5
$ws.Range('D5').Value = '''310.98'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.64%  '

# Row 6
$ws.Range('D6').Value = '''101.84'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +4.77%  '

# Row 7
$ws.Range('E7').Value = '  +1.45%  '

# Row 8
$ws.Range('E8').Value = '  -0.04%  '

# Row 9
$ws.Range('E9').Value = '  +7.61%  '

# Row 10
$ws.Range('D10').Value = '''35.72'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +1.73%  '

# Row 11
$ws.Range('E11').Value = '  +2.87%  '

# Row 12
$ws.Range('E12').Value = '  -0.70%  '

# Row 13
$ws.Range('D13').Value = '''6.99'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.90%  '

# Row 14
$ws.Range('D14').Value = '2.666.87'
$ws.Range('E14').Value = '  +1.73%  '

# Row 15
$ws.Range('D15').Value = '''14.95'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +1.22%  '

# Row 16
$ws.Range('D16').Value = '2.308.47'
$ws.Range('E16').Value = '  +1.95%  '

# Row 17
$ws.Range('E17').Value = '  +1.99%  '

# Row 18
$ws.Range('D18').Value = '43.319.92'

# Row 19
$ws.Range('D19').Value = '''12.30'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.08%  '

# Row 20
$ws.Range('D20').Value = '0.0₃0928'
$ws.Range('E20').Value = '  +2.44%  '

# Row 21
$ws.Range('E21').Value = '  +2.22%  '

# Row 22
$ws.Range('D22').Value = '''68.06'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.17%  '

# Row 23
$ws.Range('D23').Value = '''241.38'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.53%  '

# Row 24
$ws.Range('D24').Value = '''2.01'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +2.57%  '

# Row 25
$ws.Range('D25').Value = '''2.61'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +1.68%  '

# Row 27
$ws.Range('D27').Value = '''3.98'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.86%  '

# Row 28
$ws.Range('D28').Value = '''24.58'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +4.45%  '

# Row 29
$ws.Range('D29').Value = '''36.65'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -3.30%  '

# Row 30
$ws.Range('D30').Value = '''9.64'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.09%  '

# Row 31
$ws.Range('E31').Value = '  +0.39%  '

# Row 32
$ws.Range('D32').Value = '''167.43'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +3.87%  '

# Row 33
$ws.Range('E33').Value = '  +0.84%  '

# Row 35
$ws.Range('E35').Value = '  +0.77%  '

# Row 36
$ws.Range('E36').Value = '  +5.89%  '

# Row 37
$ws.Range('B37').Value = 'LidoDAOToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D37').Value = '''3.06'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -3.00%  '

# Row 38
$ws.Range('B38').Value = 'Celestia'
$ws.Range('C38').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D38').Value = '''17.58'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.38%  '

# Row 39
$ws.Range('E39').Value = '  +3.39%  '

# Row 40
$ws.Range('E40').Value = '  +1.33%  '

# Row 41
$ws.Range('E41').Value = '  +1.29%  '

# Row 42
$ws.Range('E42').Value = '  +7.22%  '

# Row 43
$ws.Range('D43').Value = '''2.32'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.66%  '

# Row 44
$ws.Range('E44').Value = '  +2.73%  '

# Row 45
$ws.Range('D45').Value = '1.970.66'
$ws.Range('E45').Value = '  +1.19%  '

# Row 46
$ws.Range('D46').Value = '''19.21'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.59%  '

# Row 47
$ws.Range('D47').Value = '''2.97'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +1.94%  '

# Row 48
$ws.Range('D48').Value = '''9.92'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.58%  '

# Row 49
$ws.Range('E49').Value = '  +3.63%  '

# Row 50
$ws.Range('D50').Value = '''2.90'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +3.14%  '

# Row 51
$ws.Range('E51').Value = '  +6.90%  '
